# Append latest Lancers scrape batch (2025-10-16 01:18 JST) to the 'ランサーズ' sheet.
# New listings are inserted in score order, pushing the lower-priority rows down,
# and every row's timestamp is refreshed to the new fetch time.
$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item("ランサーズ").Activate()
$ws = $wb.ActiveSheet

# Drop the old hyperlink objects up front; they get rebuilt below at the rows'
# new positions (plain Range writes would otherwise leave the old anchors in place).
$ws.Range("F2:F15").Hyperlinks.Delete()

# Row 2: n8n×Python×AIで公開レポート自動探索・抽出・分類・登録フロー(PoC開発)
$ws.Range("A2").Value = '2025-10-16 01:18:03'
$ws.Range("B2").Value = 'n8n×Python×AIで公開レポート自動探索・抽出・分類・登録フロー(PoC開発)'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5413825'
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5413825')
$ws.Range("G2").Value = 538
$ws.Range("H2").Value = '🔥AI,Python ◆開発'

# Row 3: 【急募】AI×LINE開発をリード!医療機関向けアプリのサーバーサイドエンジニア募集(フルリモート)
$ws.Range("A3").Value = '2025-10-16 01:18:03'
$ws.Range("B3").Value = '【急募】AI×LINE開発をリード!医療機関向けアプリのサーバーサイドエンジニア募集(フルリモート)'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5413230'
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5413230')
$ws.Range("G3").Value = 385
$ws.Range("H3").Value = '🔥AI,Ai ◆開発 ◇アプリ'

# Row 4: 【急募】不動産評価システムのAIチャットポット開発者募集
$ws.Range("A4").Value = '2025-10-16 01:18:03'
$ws.Range("B4").Value = '【急募】不動産評価システムのAIチャットポット開発者募集'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5413280'
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5413280')
$ws.Range("G4").Value = 383
$ws.Range("H4").Value = '🔥AI,Ai ◆開発'

# Row 5: 【急募】ローカルAI開発プロジェクトの協力者を探しています!
$ws.Range("A5").Value = '2025-10-16 01:18:03'
$ws.Range("B5").Value = '【急募】ローカルAI開発プロジェクトの協力者を探しています!'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5413402'
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5413402')
$ws.Range("G5").Value = 375
$ws.Range("H5").Value = '🔥AI,Ai ◆開発'

# Row 6: 【高報酬/リモート可/法人可】グローバルHRベンチャーでAIを活用し業務効率化を推進してくださる方!
$ws.Range("A6").Value = '2025-10-16 01:18:03'
$ws.Range("B6").Value = '【高報酬/リモート可/法人可】グローバルHRベンチャーでAIを活用し業務効率化を推進してくださる方!'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5413210'
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5413210')
$ws.Range("G6").Value = 370
$ws.Range("H6").Value = '🔥AI,Ai ◆効率化'

# Row 7: 【急募】AIテキスト抜粋アプリのプロンプト最適化依頼
$ws.Range("A7").Value = '2025-10-16 01:18:03'
$ws.Range("B7").Value = '【急募】AIテキスト抜粋アプリのプロンプト最適化依頼'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5413215'
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5413215')
$ws.Range("G7").Value = 338
$ws.Range("H7").Value = '🔥AI,Ai ◇アプリ'

# Row 8: 【Azure/RAG】社内文書検索AIチャットボットの精度向上&内製化支援パートナー募集!
$ws.Range("A8").Value = '2025-10-16 01:18:03'
$ws.Range("B8").Value = '【Azure/RAG】社内文書検索AIチャットボットの精度向上&内製化支援パートナー募集!'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5413954'
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5413954')
$ws.Range("G8").Value = 310
$ws.Range("H8").Value = '🔥AI,Ai'

# Row 9: 生成AIの技術顧問を募集!事業の技術選定をリードするAI専門家を募集! 【週1日〜/フルリモート】
$ws.Range("A9").Value = '2025-10-16 01:18:03'
$ws.Range("B9").Value = '生成AIの技術顧問を募集!事業の技術選定をリードするAI専門家を募集! 【週1日〜/フルリモート】'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5413955'
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5413955')
$ws.Range("G9").Value = 303
$ws.Range("H9").Value = '🔥AI,Ai'

# Row 10: コスパスポーツジムの高速自動予約botの開発
$ws.Range("A10").Value = '2025-10-16 01:18:03'
$ws.Range("B10").Value = 'コスパスポーツジムの高速自動予約botの開発'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5413835'
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5413835')
$ws.Range("G10").Value = 173
$ws.Range("H10").Value = '★bot ◆開発'

# Row 11: 3Dプリント用データのWeb自動チェック&変換&カラー補正ツール|開発パートナー募集
$ws.Range("A11").Value = '2025-10-16 01:18:03'
$ws.Range("B11").Value = '3Dプリント用データのWeb自動チェック&変換&カラー補正ツール|開発パートナー募集'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '5,000,000 円 ~ / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5413508'
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5413508')
$ws.Range("G11").Value = 135
$ws.Range("H11").Value = '◆ツール,開発'

# Row 12: セレニウムを用いた自動発注ツールの修正・機能追加
$ws.Range("A12").Value = '2025-10-16 01:18:03'
$ws.Range("B12").Value = 'セレニウムを用いた自動発注ツールの修正・機能追加'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5413916'
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5413916')
$ws.Range("G12").Value = 73
$ws.Range("H12").Value = '◆ツール'

# Row 13: IB報酬を得るための高性能EA開発依頼
$ws.Range("A13").Value = '2025-10-16 01:18:03'
$ws.Range("B13").Value = 'IB報酬を得るための高性能EA開発依頼'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5413293'
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5413293')
$ws.Range("G13").Value = 68
$ws.Range("H13").Value = '◆開発'

# Row 14: 【音声コマンド起動】超小型・低電力レコーダーのプロトタイプ開発
$ws.Range("A14").Value = '2025-10-16 01:18:03'
$ws.Range("B14").Value = '【音声コマンド起動】超小型・低電力レコーダーのプロトタイプ開発'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5413958'
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5413958')
$ws.Range("G14").Value = 63
$ws.Range("H14").Value = '◆開発'

# Row 15: wordpressレンダリングを妨げるリソースの除外
$ws.Range("A15").Value = '2025-10-16 01:18:03'
$ws.Range("B15").Value = 'wordpressレンダリングを妨げるリソースの除外'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5016989'
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5016989')
$ws.Range("G15").Value = 33
$ws.Range("H15").Value = '○WordPress'

# Row 16: 【急募】Cloud RunでWordPress構築のプロを探しています!
$ws.Range("A16").Value = '2025-10-16 01:18:03'
$ws.Range("B16").Value = '【急募】Cloud RunでWordPress構築のプロを探しています!'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Range("F16").Value = 'https://www.lancers.jp/work/detail/5413043'
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5413043')
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = '○WordPress'

# Row 17: 【急募】16タイプ診断コンテンツのLP制作
$ws.Range("A17").Value = '2025-10-16 01:18:03'
$ws.Range("B17").Value = '【急募】16タイプ診断コンテンツのLP制作'
$ws.Range("C17").Value = 'システム開発'
$ws.Range("D17").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E17").Value = '期限情報なし'
$ws.Range("F17").Value = 'https://www.lancers.jp/work/detail/5408735'
$ws.Hyperlinks.Add($ws.Range("F17"), 'https://www.lancers.jp/work/detail/5408735')
$ws.Range("G17").Value = 25

# Row 18: 【急募】Teams Roomsの設定設置と保守サポート依頼
$ws.Range("A18").Value = '2025-10-16 01:18:03'
$ws.Range("B18").Value = '【急募】Teams Roomsの設定設置と保守サポート依頼'
$ws.Range("C18").Value = 'システム開発'
$ws.Range("D18").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E18").Value = '期限情報なし'
$ws.Range("F18").Value = 'https://www.lancers.jp/work/detail/5408814'
$ws.Hyperlinks.Add($ws.Range("F18"), 'https://www.lancers.jp/work/detail/5408814')
$ws.Range("G18").Value = 18

# Row 19: Access 32bitから64bitへの修正改善依頼
$ws.Range("A19").Value = '2025-10-16 01:18:03'
$ws.Range("B19").Value = 'Access 32bitから64bitへの修正改善依頼'
$ws.Range("C19").Value = 'システム開発'
$ws.Range("D19").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E19").Value = '期限情報なし'
$ws.Range("F19").Value = 'https://www.lancers.jp/work/detail/5413333'
$ws.Hyperlinks.Add($ws.Range("F19"), 'https://www.lancers.jp/work/detail/5413333')
$ws.Range("G19").Value = 10
